$wb = $excel.ActiveWorkbook

# Column F ("想去人数") updates that apply identically to both the
# "展览" sheet and the "全部类型" sheet.
$updates = @{
    2  = 8320
    3  = 7750
    4  = 120
    9  = 116
    10 = 161
    11 = 229
    12 = 703
    14 = 1311
    19 = 118
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
